$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.162809
$ws.Range("H2").Value = 6.488427000000001
$ws.Range("I2").Value = 0.06755089002018773
$ws.Range("J2").Value = 0.06755089002018773
$ws.Range("M2").Value = 57.65261933333333
$ws.Range("N2").Value = 172.957858
$ws.Range("O2").Value = 0.6817060950001529
$ws.Range("P2").Value = 0.6817060950001529
$ws.Range("Q2").Value = 124.6916039677073
$ws.Range("R2").Value = 1122.224435709366
$ws.Range("S2").Value = 0.04604985344944697
$ws.Range("T2").Value = 0.04604985344944697

# Row 3
$ws.Range("G3").Value = 2.162809
$ws.Range("H3").Value = 6.488427000000001
$ws.Range("I3").Value = 0.06755089002018773
$ws.Range("J3").Value = 0.06755089002018773
$ws.Range("O3").Value = 0.1019529789289588
$ws.Range("P3").Value = 0.1019529789289588
$ws.Range("Q3").Value = 18.64833036579333
$ws.Range("R3").Value = 167.83497329214
$ws.Range("S3").Value = 0.006887014466860615
$ws.Range("T3").Value = 0.006887014466860615

# Row 4
$ws.Range("G4").Value = 2.162809
$ws.Range("H4").Value = 6.488427000000001
$ws.Range("I4").Value = 0.06755089002018773
$ws.Range("J4").Value = 0.06755089002018773
$ws.Range("M4").Value = 2.790736
$ws.Range("N4").Value = 8.372208
$ws.Range("O4").Value = 0.0329987043561157
$ws.Range("P4").Value = 0.0329987043561157
$ws.Range("Q4").Value = 6.035828937424001
$ws.Range("R4").Value = 54.32246043681601
$ws.Range("S4").Value = 0.002229091848768662
$ws.Range("T4").Value = 0.002229091848768661

# Row 5
$ws.Range("G5").Value = 2.162809
$ws.Range("H5").Value = 6.488427000000001
$ws.Range("I5").Value = 0.06755089002018773
$ws.Range("J5").Value = 0.06755089002018773
$ws.Range("M5").Value = 15.50544933333333
$ws.Range("N5").Value = 46.516348
$ws.Range("O5").Value = 0.1833422217147727
$ws.Range("P5").Value = 0.1833422217147727
$ws.Range("Q5").Value = 33.53532536717734
$ws.Range("R5").Value = 301.817928304596
$ws.Range("S5").Value = 0.01238493025511148
$ws.Range("T5").Value = 0.01238493025511148

# Row 6
$ws.Range("I6").Value = 0.5628021396814664
$ws.Range("J6").Value = 0.5628021396814664
$ws.Range("M6").Value = 57.65261933333333
$ws.Range("N6").Value = 172.957858
$ws.Range("O6").Value = 0.6817060950001529
$ws.Range("P6").Value = 0.6817060950001529
$ws.Range("Q6").Value = 1038.871604687477
$ws.Range("R6").Value = 9349.844442187295
$ws.Range("S6").Value = 0.3836656488999831
$ws.Range("T6").Value = 0.3836656488999831

# Row 7
$ws.Range("I7").Value = 0.5628021396814664
$ws.Range("J7").Value = 0.5628021396814664
$ws.Range("O7").Value = 0.1019529789289588
$ws.Range("P7").Value = 0.1019529789289588
$ws.Range("S7").Value = 0.05737935468811749
$ws.Range("T7").Value = 0.05737935468811749

# Row 8
$ws.Range("I8").Value = 0.5628021396814664
$ws.Range("J8").Value = 0.5628021396814664
$ws.Range("M8").Value = 2.790736
$ws.Range("N8").Value = 8.372208
$ws.Range("O8").Value = 0.0329987043561157
$ws.Range("P8").Value = 0.0329987043561157
$ws.Range("Q8").Value = 50.28767851494401
$ws.Range("R8").Value = 452.589106634496
$ws.Range("S8").Value = 0.01857174141833804
$ws.Range("T8").Value = 0.01857174141833804

# Row 9
$ws.Range("I9").Value = 0.5628021396814664
$ws.Range("J9").Value = 0.5628021396814664
$ws.Range("M9").Value = 15.50544933333333
$ws.Range("N9").Value = 46.516348
$ws.Range("O9").Value = 0.1833422217147727
$ws.Range("P9").Value = 0.1833422217147727
$ws.Range("Q9").Value = 279.4005062837974
$ws.Range("R9").Value = 2514.604556554176
$ws.Range("S9").Value = 0.1031853946750279
$ws.Range("T9").Value = 0.1031853946750279

# Row 10
$ws.Range("G10").Value = 4.650307000000001
$ws.Range("H10").Value = 13.950921
$ws.Range("I10").Value = 0.1452427730405732
$ws.Range("J10").Value = 0.1452427730405732
$ws.Range("M10").Value = 57.65261933333333
$ws.Range("N10").Value = 172.957858
$ws.Range("O10").Value = 0.6817060950001529
$ws.Range("P10").Value = 0.6817060950001529
$ws.Range("Q10").Value = 268.1023792541353
$ws.Range("R10").Value = 2412.921413287218
$ws.Range("S10").Value = 0.09901288363648265
$ws.Range("T10").Value = 0.09901288363648265

# Row 11
$ws.Range("G11").Value = 4.650307000000001
$ws.Range("H11").Value = 13.950921
$ws.Range("I11").Value = 0.1452427730405732
$ws.Range("J11").Value = 0.1452427730405732
$ws.Range("O11").Value = 0.1019529789289588
$ws.Range("P11").Value = 0.1019529789289588
$ws.Range("Q11").Value = 40.09621803791333
$ws.Range("R11").Value = 360.86596234122
$ws.Range("S11").Value = 0.01480793337938911
$ws.Range("T11").Value = 0.01480793337938911

# Row 12
$ws.Range("G12").Value = 4.650307000000001
$ws.Range("H12").Value = 13.950921
$ws.Range("I12").Value = 0.1452427730405732
$ws.Range("J12").Value = 0.1452427730405732
$ws.Range("M12").Value = 2.790736
$ws.Range("N12").Value = 8.372208
$ws.Range("O12").Value = 0.0329987043561157
$ws.Range("P12").Value = 0.0329987043561157
$ws.Range("Q12").Value = 12.977779155952
$ws.Range("R12").Value = 116.800012403568
$ws.Range("S12").Value = 0.004792823327428288
$ws.Range("T12").Value = 0.004792823327428287

# Row 13
$ws.Range("G13").Value = 4.650307000000001
$ws.Range("H13").Value = 13.950921
$ws.Range("I13").Value = 0.1452427730405732
$ws.Range("J13").Value = 0.1452427730405732
$ws.Range("M13").Value = 15.50544933333333
$ws.Range("N13").Value = 46.516348
$ws.Range("O13").Value = 0.1833422217147727
$ws.Range("P13").Value = 0.1833422217147727
$ws.Range("Q13").Value = 72.10509957294533
$ws.Range("R13").Value = 648.9458961565081
$ws.Range("S13").Value = 0.02662913269727319
$ws.Range("T13").Value = 0.02662913269727319

# Row 14
$ws.Range("G14").Value = 7.184856000000001
$ws.Range("H14").Value = 21.554568
$ws.Range("I14").Value = 0.2244041972577726
$ws.Range("J14").Value = 0.2244041972577726
$ws.Range("M14").Value = 57.65261933333333
$ws.Range("N14").Value = 172.957858
$ws.Range("O14").Value = 0.6817060950001529
$ws.Range("P14").Value = 0.6817060950001529
$ws.Range("Q14").Value = 414.225767932816
$ws.Range("R14").Value = 3728.031911395344
$ws.Range("S14").Value = 0.1529777090142402
$ws.Range("T14").Value = 0.1529777090142402

# Row 15
$ws.Range("G15").Value = 7.184856000000001
$ws.Range("H15").Value = 21.554568
$ws.Range("I15").Value = 0.2244041972577726
$ws.Range("J15").Value = 0.2244041972577726
$ws.Range("O15").Value = 0.1019529789289588
$ws.Range("P15").Value = 0.1019529789289588
$ws.Range("Q15").Value = 61.94979229264
$ws.Range("R15").Value = 557.54813063376
$ws.Range("S15").Value = 0.02287867639459161
$ws.Range("T15").Value = 0.02287867639459161

# Row 16
$ws.Range("G16").Value = 7.184856000000001
$ws.Range("H16").Value = 21.554568
$ws.Range("I16").Value = 0.2244041972577726
$ws.Range("J16").Value = 0.2244041972577726
$ws.Range("M16").Value = 2.790736
$ws.Range("N16").Value = 8.372208
$ws.Range("O16").Value = 0.0329987043561157
$ws.Range("P16").Value = 0.0329987043561157
$ws.Range("Q16").Value = 20.051036294016
$ws.Range("R16").Value = 180.459326646144
$ws.Range("S16").Value = 0.007405047761580709
$ws.Range("T16").Value = 0.007405047761580707

# Row 17
$ws.Range("G17").Value = 7.184856000000001
$ws.Range("H17").Value = 21.554568
$ws.Range("I17").Value = 0.2244041972577726
$ws.Range("J17").Value = 0.2244041972577726
$ws.Range("M17").Value = 15.50544933333333
$ws.Range("N17").Value = 46.516348
$ws.Range("O17").Value = 0.1833422217147727
$ws.Range("P17").Value = 0.1833422217147727
$ws.Range("Q17").Value = 111.404420675296
$ws.Range("R17").Value = 1002.639786077664
$ws.Range("S17").Value = 0.04114276408736013
$ws.Range("T17").Value = 0.04114276408736013
